$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increase MaxInvest Storage Adapt Szenarios Existing Units (column F = ExisUnits)
$ws.Range("F8").Value = 17
$ws.Range("F10").Value = 1
$ws.Range("F11").Value = 2
$ws.Range("F12").Value = 16
$ws.Range("F14").Value = 6
$ws.Range("F15").Value = 13
$ws.Range("F16").Value = 78

$ws.Range("I21").Select()
